$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows whose target cluster is "Resolving-Mac" (rows 5, 9, 13
# in the original 1-header+12-data-row table). Delete from the bottom up so the
# remaining row numbers do not shift while we are still removing earlier ones.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(5).Delete()

# Refresh the 9 remaining data rows (now rows 2-10) with the recomputed TPM-based
# ligand/receptor/edge expression values. Sending/target cluster labels, ligand,
# receptor and the E/F/K/L count columns are unchanged by this update.

$ws.Range("G2").Value = 6.670167666666667
$ws.Range("H2").Value = 20.010503
$ws.Range("I2").Value = 0.0423069620011633
$ws.Range("J2").Value = 0.0423069620011633
$ws.Range("M2").Value = 0.274713
$ws.Range("N2").Value = 0.824139
$ws.Range("O2").Value = 0.1055967877339779
$ws.Range("P2").Value = 0.1055967877339779
$ws.Range("Q2").Value = 1.832381770213
$ws.Range("R2").Value = 16.491435931917
$ws.Range("S2").Value = 0.004467479286106309
$ws.Range("T2").Value = 0.004467479286106309

$ws.Range("G3").Value = 6.670167666666667
$ws.Range("H3").Value = 20.010503
$ws.Range("I3").Value = 0.0423069620011633
$ws.Range("J3").Value = 0.0423069620011633
$ws.Range("M3").Value = 0.8886716666666666
$ws.Range("N3").Value = 2.666015
$ws.Range("O3").Value = 0.3415960415058637
$ws.Range("P3").Value = 0.3415960415058638
$ws.Range("Q3").Value = 5.927589017282778
$ws.Range("R3").Value = 53.34830115554499
$ws.Range("S3").Value = 0.01445189074773637
$ws.Range("T3").Value = 0.01445189074773638

$ws.Range("G4").Value = 6.670167666666667
$ws.Range("H4").Value = 20.010503
$ws.Range("I4").Value = 0.0423069620011633
$ws.Range("J4").Value = 0.0423069620011633
$ws.Range("M4").Value = 1.438143333333333
$ws.Range("N4").Value = 4.31443
$ws.Range("O4").Value = 0.5528071707601584
$ws.Range("P4").Value = 0.5528071707601584
$ws.Range("Q4").Value = 9.592657162032223
$ws.Range("R4").Value = 86.33391445829
$ws.Range("S4").Value = 0.02338759196732061
$ws.Range("T4").Value = 0.02338759196732061

$ws.Range("G5").Value = 149.9875183333334
$ws.Range("H5").Value = 449.9625550000001
$ws.Range("I5").Value = 0.9513278459982415
$ws.Range("J5").Value = 0.9513278459982416
$ws.Range("M5").Value = 0.274713
$ws.Range("N5").Value = 0.824139
$ws.Range("O5").Value = 0.1055967877339779
$ws.Range("P5").Value = 0.1055967877339779
$ws.Range("Q5").Value = 41.20352112390501
$ws.Range("R5").Value = 370.831690115145
$ws.Range("S5").Value = 0.1004571646192987
$ws.Range("T5").Value = 0.1004571646192987

$ws.Range("G6").Value = 149.9875183333334
$ws.Range("H6").Value = 449.9625550000001
$ws.Range("I6").Value = 0.9513278459982415
$ws.Range("J6").Value = 0.9513278459982416
$ws.Range("M6").Value = 0.8886716666666666
$ws.Range("N6").Value = 2.666015
$ws.Range("O6").Value = 0.3415960415058637
$ws.Range("P6").Value = 0.3415960415058638
$ws.Range("Q6").Value = 133.2896578964806
$ws.Range("R6").Value = 1199.606921068325
$ws.Range("S6").Value = 0.3249698263672992
$ws.Range("T6").Value = 0.3249698263672993

$ws.Range("G7").Value = 149.9875183333334
$ws.Range("H7").Value = 449.9625550000001
$ws.Range("I7").Value = 0.9513278459982415
$ws.Range("J7").Value = 0.9513278459982416
$ws.Range("M7").Value = 1.438143333333333
$ws.Range("N7").Value = 4.31443
$ws.Range("O7").Value = 0.5528071707601584
$ws.Range("P7").Value = 0.5528071707601584
$ws.Range("Q7").Value = 215.7035495742945
$ws.Range("R7").Value = 1941.33194616865
$ws.Range("S7").Value = 0.5259008550116435
$ws.Range("T7").Value = 0.5259008550116436

$ws.Range("G8").Value = 1.003544
$ws.Range("H8").Value = 3.010632
$ws.Range("I8").Value = 0.0063651920005952
$ws.Range("J8").Value = 0.006365192000595201
$ws.Range("M8").Value = 0.274713
$ws.Range("N8").Value = 0.824139
$ws.Range("O8").Value = 0.1055967877339779
$ws.Range("P8").Value = 0.1055967877339779
$ws.Range("Q8").Value = 0.275686582872
$ws.Range("R8").Value = 2.481179245848
$ws.Range("S8").Value = 0.0006721438285728653
$ws.Range("T8").Value = 0.0006721438285728654

$ws.Range("G9").Value = 1.003544
$ws.Range("H9").Value = 3.010632
$ws.Range("I9").Value = 0.0063651920005952
$ws.Range("J9").Value = 0.006365192000595201
$ws.Range("M9").Value = 0.8886716666666666
$ws.Range("N9").Value = 2.666015
$ws.Range("O9").Value = 0.3415960415058637
$ws.Range("P9").Value = 0.3415960415058638
$ws.Range("Q9").Value = 0.8918211190533333
$ws.Range("R9").Value = 8.02639007148
$ws.Range("S9").Value = 0.00217432439082811
$ws.Range("T9").Value = 0.00217432439082811

$ws.Range("G10").Value = 1.003544
$ws.Range("H10").Value = 3.010632
$ws.Range("I10").Value = 0.0063651920005952
$ws.Range("J10").Value = 0.0063651920005952
$ws.Range("M10").Value = 1.438143333333333
$ws.Range("N10").Value = 4.31443
$ws.Range("O10").Value = 0.5528071707601584
$ws.Range("P10").Value = 0.5528071707601584
$ws.Range("Q10").Value = 1.443240113306667
$ws.Range("R10").Value = 12.98916101976
$ws.Range("S10").Value = 0.003518723781194225
$ws.Range("T10").Value = 0.003518723781194225
